# Apply "update vals with new BTT BCT grid" changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the BTT -> BCT (column D) data values with the new grid counts.
$ws.Range("D4").Value = 21
$ws.Range("D6").Value = 11
$ws.Range("D7").Value = 4
$ws.Range("D8").Value = 3

# Update the view: zoom in on the sheet and move the active selection to D10.
$ws.Activate()
$excel.ActiveWindow.Zoom = 160
$ws.Range("D10").Select()
